$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D19").Value = "controlli sulla scrittura dei numeri/liste e MessageDialog fine gioco"

$ws.Range("D23").ClearContents()
$ws.Range("F23").ClearContents()
$ws.Range("J24").ClearContents()
$ws.Range("K24").ClearContents()

$ws.Columns.Item(4).ColumnWidth = 59.5

$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("F27").Select()
